$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = " Abu Dhabi"
$ws.Range("B11").Value = " October 28 2020"
$ws.Range("C11").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D11").Value = "Royal Challengers Bangalore"
$ws.Range("E11").Value = "Mumbai Indians"
$ws.Range("F11").Value = "Shivam Dube "
$ws.Range("G11").Value = "'2"
$ws.Range("H11").Value = "'6"
$ws.Range("I11").Value = "'0"
$ws.Range("J11").Value = "'0"
$ws.Range("K11").Value = "'33.33"
$ws.Range("G11:K11").Style = "Normal"

$ws.Range("A12").Value = " Abu Dhabi"
$ws.Range("B12").Value = " November 02 2020"
$ws.Range("C12").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D12").Value = "Royal Challengers Bangalore"
$ws.Range("E12").Value = "Delhi Capitals"
$ws.Range("F12").Value = "Shivam Dube "
$ws.Range("G12").Value = "'17"
$ws.Range("H12").Value = "'11"
$ws.Range("I12").Value = "'2"
$ws.Range("J12").Value = "'1"
$ws.Range("K12").Value = "'154.54"
$ws.Range("G12:K12").Style = "Normal"

$ws.Range("A13").Value = " Dubai (DSC)"
$ws.Range("B13").Value = " September 24 2020"
$ws.Range("C13").Value = "Kings XI won by 97 runs"
$ws.Range("D13").Value = "Royal Challengers Bangalore"
$ws.Range("E13").Value = "Kings XI Punjab"
$ws.Range("F13").Value = "Shivam Dube "
$ws.Range("G13").Value = "'12"
$ws.Range("H13").Value = "'12"
$ws.Range("I13").Value = "'0"
$ws.Range("J13").Value = "'1"
$ws.Range("K13").Value = "'100.00"
$ws.Range("G13:K13").Style = "Normal"

$ws.Range("A14").Value = " Abu Dhabi"
$ws.Range("B14").Value = " November 06 2020"
$ws.Range("C14").Value = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$ws.Range("D14").Value = "Royal Challengers Bangalore"
$ws.Range("E14").Value = "Sunrisers Hyderabad"
$ws.Range("F14").Value = "Shivam Dube "
$ws.Range("G14").Value = "'8"
$ws.Range("H14").Value = "'13"
$ws.Range("I14").Value = "'0"
$ws.Range("J14").Value = "'0"
$ws.Range("K14").Value = "'61.53"
$ws.Range("G14:K14").Style = "Normal"

$ws.Range("A15").Value = " Sharjah"
$ws.Range("B15").Value = " October 15 2020"
$ws.Range("C15").Value = "Kings XI won by 8 wickets"
$ws.Range("D15").Value = "Royal Challengers Bangalore"
$ws.Range("E15").Value = "Kings XI Punjab"
$ws.Range("F15").Value = "Shivam Dube "
$ws.Range("G15").Value = "'23"
$ws.Range("H15").Value = "'19"
$ws.Range("I15").Value = "'0"
$ws.Range("J15").Value = "'2"
$ws.Range("K15").Value = "'121.05"
$ws.Range("G15:K15").Style = "Normal"

$ws.Range("A16").Value = " Dubai (DSC)"
$ws.Range("B16").Value = " October 05 2020"
$ws.Range("C16").Value = "Capitals won by 59 runs"
$ws.Range("D16").Value = "Royal Challengers Bangalore"
$ws.Range("E16").Value = "Delhi Capitals"
$ws.Range("F16").Value = "Shivam Dube "
$ws.Range("G16").Value = "'11"
$ws.Range("H16").Value = "'12"
$ws.Range("I16").Value = "'0"
$ws.Range("J16").Value = "'1"
$ws.Range("K16").Value = "'91.66"
$ws.Range("G16:K16").Style = "Normal"

$ws.Range("A17").Value = " Dubai (DSC)"
$ws.Range("B17").Value = " September 28 2020"
$ws.Range("C17").Value = "Match tied (RCB won the one-over eliminator)"
$ws.Range("D17").Value = "Royal Challengers Bangalore"
$ws.Range("E17").Value = "Mumbai Indians"
$ws.Range("F17").Value = "Shivam Dube "
$ws.Range("G17").Value = "'27"
$ws.Range("H17").Value = "'10"
$ws.Range("I17").Value = "'1"
$ws.Range("J17").Value = "'3"
$ws.Range("K17").Value = "'270.00"
$ws.Range("G17:K17").Style = "Normal"

$ws.Range("A18").Value = " Dubai (DSC)"
$ws.Range("B18").Value = " September 21 2020"
$ws.Range("C18").Value = "RCB won by 10 runs"
$ws.Range("D18").Value = "Royal Challengers Bangalore"
$ws.Range("E18").Value = "Sunrisers Hyderabad"
$ws.Range("F18").Value = "Shivam Dube "
$ws.Range("G18").Value = "'7"
$ws.Range("H18").Value = "'8"
$ws.Range("I18").Value = "'0"
$ws.Range("J18").Value = "'0"
$ws.Range("K18").Value = "'87.50"
$ws.Range("G18:K18").Style = "Normal"

$ws.Range("A19").Value = " Dubai (DSC)"
$ws.Range("B19").Value = " October 10 2020"
$ws.Range("C19").Value = "RCB won by 37 runs"
$ws.Range("D19").Value = "Royal Challengers Bangalore"
$ws.Range("E19").Value = "Chennai Super Kings"
$ws.Range("F19").Value = "Shivam Dube "
$ws.Range("G19").Value = "'22"
$ws.Range("H19").Value = "'14"
$ws.Range("I19").Value = "'2"
$ws.Range("J19").Value = "'1"
$ws.Range("K19").Value = "'157.14"
$ws.Range("G19:K19").Style = "Normal"

